# Adjust Investment Summary table column widths for better formatting
#
# Widens the last grid column of each of the three comparison / summary
# tables (slides 2, 3, 4) by a single EMU and clears the placeholder
# sample content (text + per-cell shading) that was left over from the
# template, leaving empty cells ready for real data.

$p = $ppt.ActivePresentation

# EMU -> point helper (PowerPoint COM table/column widths are expressed in points).
function EmuToPt([double]$emu) {
    return $emu / 12700.0
}

function Clear-TableCells($table, [int]$rows, [int]$cols) {
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $table.Cell($r, $c).Shape.TextFrame.TextRange.Text = ""
        }
    }
}

# --- Slide 2: "Why This Solution?" (Traditional Approach / Our Solution) ---
$slide2 = $p.Slides.Item(2)
$tbl2 = $slide2.Shapes.Item(3).Table
$tbl2.Columns.Item(2).Width = EmuToPt 4355467
Clear-TableCells $tbl2 4 2

# --- Slide 3: "Business Value - Financial Impact" (Metric / Value) ---
$slide3 = $p.Slides.Item(3)
$tbl3 = $slide3.Shapes.Item(3).Table
$tbl3.Columns.Item(2).Width = EmuToPt 4355467
Clear-TableCells $tbl3 6 2

# --- Slide 4: "Risk Mitigation" (Risk / Mitigation Strategy / Success Probability) ---
$slide4 = $p.Slides.Item(4)
$tbl4 = $slide4.Shapes.Item(3).Table
$tbl4.Columns.Item(3).Width = EmuToPt 2903645
Clear-TableCells $tbl4 4 3
